$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Workbook default ("Normal") font: Calibri -> 宋体 (SimSun)
# ---------------------------------------------------------------------
$wb.Styles("Normal").Font.Name = "宋体"

# ---------------------------------------------------------------------
# 2. New column C with English translations.  Two of the values start
#    with a literal apostrophe; entering that through .Value would flip
#    on Excel's "quote prefix" (text-forced) cell flag, so instead we
#    build the text with a CHAR(39) formula and flatten it down to a
#    plain value via copy / paste-special, exactly as a user would.
# ---------------------------------------------------------------------
$ws.Range("C1").Formula = "=CHAR(39)&""Please Select""&CHAR(39)"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C4").Formula = "=CHAR(39)&""All Subnets""&CHAR(39)"
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("C2").Value = "No available subnet--'"
$ws.Range("C3").Value = "No available subnet, please create a new subnet"

# ---------------------------------------------------------------------
# 3. Rich-text runs inside column B that mix 宋体 (CJK glyphs) with
#    Arial (the ASCII quote marks) -- mirrors the original authoring.
#    (Characters() calls are split so a full-length run still keeps its
#    explicit <rPr>, instead of collapsing into a cell-level style.)
# ---------------------------------------------------------------------
# B3 = 无可用子网，请新建子网  (single run, all 宋体)
$ws.Range("B3").Characters(1, 10).Font.Name = "宋体"
$ws.Range("B3").Characters(11, 1).Font.Name = "宋体"

# B1 = '请选择'
$ws.Range("B1").Characters(2, 3).Font.Name = "宋体"
$ws.Range("B1").Characters(5, 1).Font.Name = "Arial"

# B2 = 无可用子网--'
$ws.Range("B2").Characters(1, 5).Font.Name = "宋体"
$ws.Range("B2").Characters(6, 3).Font.Name = "Arial"

# B4 = '全部子网'
$ws.Range("B4").Characters(2, 4).Font.Name = "宋体"
$ws.Range("B4").Characters(6, 1).Font.Name = "Arial"

# ---------------------------------------------------------------------
# 4. Whole data range uses Arial as its cell-level font
# ---------------------------------------------------------------------
$ws.Range("A1:C4").Font.Name = "Arial"

# ---------------------------------------------------------------------
# 5. Column widths (approximate -- Excel quantises to pixel widths)
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 14.1667
$ws.Columns("B").ColumnWidth = 24.1667
$ws.Columns("C").ColumnWidth = 51.9167

# ---------------------------------------------------------------------
# 6. Sheet view / selection cosmetics
# ---------------------------------------------------------------------
$ws.Range("C19").Select() | Out-Null
